# Applies the cryptos-list price/volume refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference -> new text value.
# Values that Excel would auto-parse as a plain number (e.g. "0.501") are written
# with the cell temporarily forced to Text format so they stay strings (matching the
# source "0.500" -> "0.501" style text cells), then the number format is cleared again
# so the cell style is left exactly as it was before (General / no explicit format).
$updates = [ordered]@{
    'D2' = '60.346.36'
    'E2' = '  -5.88%  '
    'D3' = '3.009.57'
    'E3' = '  -6.26%  '
    'E4' = '  +0.10%  '
    'D5' = '573.06'
    'E5' = '  -4.01%  '
    'D6' = '126.92'
    'E6' = '  -8.05%  '
    'E7' = '  +0.06%  '
    'D8' = '3.005.84'
    'E8' = '  -6.29%  '
    'D9' = '0.501'
    'E9' = '  -2.79%  '
    'D10' = '0.131'
    'E10' = '  -9.07%  '
    'D11' = '5.13'
    'E11' = '  -4.49%  '
    'E12' = '  -4.13%  '
    'D13' = '0.0000219'
    'E13' = '  -9.57%  '
    'D14' = '32.62'
    'E14' = '  -7.01%  '
    'E15' = '  +0.33%  '
    'D16' = '3.509.00'
    'E16' = '  -6.27%  '
    'D17' = '3.010.54'
    'E17' = '  -6.32%  '
    'D18' = '60.331.01'
    'E18' = '  -5.84%  '
    'D19' = '6.41'
    'E19' = '  -2.71%  '
    'D20' = '429.02'
    'E20' = '  -7.71%  '
    'D21' = '13.12'
    'E21' = '  -6.62%  '
    'D22' = '0.668'
    'E22' = '  -4.93%  '
    'D23' = '7.04'
    'E23' = '  -8.49%  '
    'D24' = '13.06'
    'E24' = '  -1.60%  '
    'D25' = '79.27'
    'E25' = '  -5.26%  '
    'E26' = '  +0.11%  '
    'D27' = '0.998'
    'E27' = '  -0.16%  '
    'D28' = '2.55'
    'E28' = '  -5.63%  '
    'E29' = '  -7.57%  '
    'D30' = '1.94'
    'E30' = '  -8.14%  '
    'B31' = 'NEARProtocol'
    'C31' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    'D31' = '6.13'
    'E31' = '  -10.29%  '
    'B32' = 'EthereumClassic'
    'C32' = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
    'D32' = '25.21'
    'E32' = '  -9.22%  '
    'D33' = '0.0935'
    'E33' = '  -8.61%  '
    'D34' = '0.950'
    'E34' = '  -8.10%  '
    'D35' = '5.62'
    'E35' = '  -5.16%  '
    'B36' = 'Stacks'
    'C36' = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
    'D36' = '2.08'
    'E36' = '  -16.72%  '
    'B37' = 'OKB'
    'C37' = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
    'D37' = '50.38'
    'E37' = '  -2.69%  '
    'B38' = 'PEPE'
    'C38' = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
    'D38' = '0.0₃0663'
    'E38' = '  -11.38%  '
    'B39' = 'Cosmos'
    'C39' = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
    'D39' = '8.30'
    'E39' = '  +1.82%  '
    'D40' = '387.34'
    'E40' = '  -3.23%  '
    'D41' = '0.0355'
    'E41' = '  -9.97%  '
    'D42' = '0.110'
    'E42' = '  -3.42%  '
    'D43' = '2.52'
    'E43' = '  -8.21%  '
    'D44' = '2.661.50'
    'E44' = '  -5.75%  '
    'E45' = '  +0.06%  '
    'B46' = 'TheGraph'
    'C46' = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
    'D46' = '0.235'
    'E46' = '  -8.09%  '
    'B47' = 'Fetch.AI'
    'C47' = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
    'D47' = '2.01'
    'E47' = '  -7.23%  '
    'D48' = '120.03'
    'E48' = '  -8.51%  '
    'E49' = '  -4.16%  '
    'D50' = '23.66'
    'E50' = '  -8.21%  '
    'D51' = '0.133'
    'E51' = '  +0.39%  '
}

foreach ($cellRef in $updates.Keys) {
    $value = $updates[$cellRef]
    $cell = $ws.Range($cellRef)
    $numericLooking = $value -match '^[+-]?[0-9]*\.?[0-9]+$'
    if ($numericLooking) {
        # Force text so the numeric-looking string is not reinterpreted as a number.
        $cell.NumberFormat = "@"
        $cell.Value = $value
        # Restore the default (General) number format so no stray style is left behind.
        $cell.ClearFormats()
    } else {
        $cell.Value = $value
    }
}

Write-Output "Applied $($updates.Count) cell updates"
